$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "70.964.83"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.572.16"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "584.10"
$ws.Range("E5").Value = "  +2.51%  "
Set-TextValue "D6" "186.48"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("D7").Value = "3.561.69"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +16.63%  "
Set-TextValue "D11" "0.654"
$ws.Range("E11").Value = "  +2.67%  "
Set-TextValue "D12" "54.76"
$ws.Range("E12").Value = "  +1.76%  "
Set-TextValue "D13" "0.0000320"
$ws.Range("E13").Value = "  +6.31%  "
Set-TextValue "D14" "9.51"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "4.130.54"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "70.933.40"
$ws.Range("E16").Value = "  +2.71%  "
Set-TextValue "D17" "19.34"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "3.529.68"
$ws.Range("E18").Value = "  +0.77%  "
Set-TextValue "D19" "12.45"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D20" "0.121"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "561.76"
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("E23").Value = "  -13.49%  "
Set-TextValue "D24" "5.05"
$ws.Range("E24").Value = "  +1.47%  "
Set-TextValue "D25" "4.59"
$ws.Range("E25").Value = "  +5.10%  "
Set-TextValue "D26" "94.22"
$ws.Range("E26").Value = "  +0.48%  "
Set-TextValue "D27" "11.33"
$ws.Range("E27").Value = "  +2.82%  "
Set-TextValue "D28" "2.96"
Set-TextValue "D29" "9.17"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  +2.89%  "
Set-TextValue "D31" "7.30"
$ws.Range("E31").Value = "  +0.51%  "
Set-TextValue "D32" "12.34"
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("E33").Value = "  +2.93%  "
Set-TextValue "D34" "63.65"
$ws.Range("E34").Value = "  -0.94%  "
Set-TextValue "D35" "3.40"
$ws.Range("E35").Value = "  +10.84%  "
Set-TextValue "D36" "554.70"
$ws.Range("E36").Value = "  -2.69%  "
Set-TextValue "D37" "0.420"
$ws.Range("E37").Value = "  +5.56%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D38" "37.89"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0804"
$ws.Range("E39").Value = "  +5.64%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +9.73%  "
$ws.Range("D42").Value = "3.562.98"
$ws.Range("E42").Value = "  +11.86%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D43" "0.137"
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "3.45"
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("E45").Value = "  +1.55%  "
Set-TextValue "D46" "3.51"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  -0.69%  "
Set-TextValue "D48" "9.38"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  +2.86%  "
Set-TextValue "D50" "1.50"
$ws.Range("E50").Value = "  +10.80%  "
Set-TextValue "D51" "0.997"
$ws.Range("E51").Value = "  +0.04%  "
